# Apply updated crypto price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.929.31'
$ws.Range("E2").Value = '  +2.03%  '
$ws.Range("D3").Value = '2.051.25'
$ws.Range("E3").Value = '  +1.21%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.80'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.52%  '
$ws.Range("E9").Value = '  +1.84%  '
$ws.Range("E10").Value = '  +3.04%  '
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").Value = '2.355.95'
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.754'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.75%  '
$ws.Range("D17").Value = '2.057.07'
$ws.Range("E17").Value = '  +3.00%  '
$ws.Range("D18").Value = '37.842.16'
$ws.Range("E18").Value = '  +2.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("D21").Value = '0.0₃0836'
$ws.Range("E21").Value = '  +2.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.39%  '
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.133'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.96'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.16%  '
$ws.Range("E31").Value = '  +1.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.59'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.21%  '
$ws.Range("E34").Value = '  +10.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0609'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.64%  '
$ws.Range("E38").Value = '  +6.12%  '
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").Value = '1.493.59'
$ws.Range("E40").Value = '  +1.94%  '
$ws.Range("E41").Value = '  +1.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.63%  '
$ws.Range("E43").Value = '  +2.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("E47").Value = '  +13.43%  '
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.28%  '
$ws.Range("D51").Value = '2.245.59'
$ws.Range("E51").Value = '  +1.38%  '
